$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.226.01'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '1.632.47'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.40'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('E6').Value = '  +1.14%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0627'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.35'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0851'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').Value = '1.637.59'
$ws.Range('E12').Value = '  -1.00%  '
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.545'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '65.18'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.78%  '
$ws.Range('D16').Value = '27.206.22'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').Value = '0.0₃0741'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '218.99'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.61%  '
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.96'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.56%  '
$ws.Range('E21').Value = '  -1.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.45'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.08'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '147.86'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.27'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.75%  '
$ws.Range('E28').Value = '  -1.19%  '
$ws.Range('E29').Value = '  -0.65%  '
$ws.Range('E30').Value = '  -0.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.39'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.26%  '
$ws.Range('D33').Value = '1.320.41'
$ws.Range('E33').Value = '  +4.30%  '
$ws.Range('E34').Value = '  -1.00%  '
$ws.Range('E35').Value = '  -0.44%  '
$ws.Range('E36').Value = '  -1.57%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.544'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.850'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E40').Value = '  +1.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.800'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '64.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.97%  '
$ws.Range('D43').Value = '1.769.40'
$ws.Range('E43').Value = '  -1.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.22'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.82%  '
$ws.Range('E46').Value = '  +0.36%  '
$ws.Range('D47').Value = '0.0₆0106'
$ws.Range('E47').Value = '  +0.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.810'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +21.61%  '
$ws.Range('E49').Value = '  +0.20%  '
$ws.Range('E50').Value = '  -0.86%  '
$ws.Range('E51').Value = '  -1.72%  '
